$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New USB connector part (J1): Amphenol ICC (FCI) 10118193-0001LF
# replaces the previous Molex 1051640001 part, with updated DigiKey part number.
$ws.Range("D7").Value = "Amphenol ICC (FCI)"
$ws.Range("E7").Value = "10118193-0001LF"
$ws.Range("I7").Value = "DigiKey Part: 609-4616-2-ND"
